$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("E").Insert()
$ws.Range("E1").Value = "logG_2025"
$ws.Range("E2:E171").Formula = "=LN(D2)"
$ws.Range("E74").ClearContents()
$ws.Range("E103").ClearContents()
$ws.Range("E109").ClearContents()
$ws.Range("E158").ClearContents()
$ws.Range("E164").ClearContents()
$ws.Range("E169").ClearContents()
$ws.Range("E170").ClearContents()
$ws.Range("E171").ClearContents()
$ws.Range("C170").ClearContents()
$ws.Range("D170").ClearContents()
$ws.Range("C171").ClearContents()
$ws.Range("D171").ClearContents()
